$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.51175556330869
$ws.Range("C2").Value = 10.86479052533408
$ws.Range("D2").Value = 9.933828302217213
$ws.Range("F2").Value = 29.1252241856137
$ws.Range("G2").Value = 28.01194080705493
$ws.Range("H2").Value = 14.20942627117649
$ws.Range("I2").Value = 21.88902293165004
$ws.Range("J2").Value = 10.19938152626566
$ws.Range("L2").Value = 11.88597945862589
$ws.Range("O2").Value = 21.50947433167648
$ws.Range("B3").Value = 16.82198735000601
$ws.Range("C3").Value = 10.61755766320685
$ws.Range("D3").Value = 9.9331822252972
$ws.Range("F3").Value = 29.25186847755267
$ws.Range("G3").Value = 28.16618813898443
$ws.Range("H3").Value = 14.27198728391407
$ws.Range("I3").Value = 22.03861408425451
$ws.Range("J3").Value = 10.23114005742055
$ws.Range("L3").Value = 11.84998952917197
$ws.Range("O3").Value = 21.62056097327883
$ws.Range("B4").Value = 16.38337600748879
$ws.Range("C4").Value = 10.46270352913917
$ws.Range("D4").Value = 9.934064610657716
$ws.Range("F4").Value = 29.33743223727568
$ws.Range("G4").Value = 28.27174015364201
$ws.Range("H4").Value = 14.31299661071504
$ws.Range("I4").Value = 22.13543576999547
$ws.Range("J4").Value = 10.25165822499705
$ws.Range("L4").Value = 11.82911104097973
$ws.Range("O4").Value = 21.69413583021914
$ws.Range("B5").Value = 16.20106789366497
$ws.Range("C5").Value = 10.3989038485347
$ws.Range("D5").Value = 9.934746935131296
$ws.Range("F5").Value = 29.37425672590074
$ws.Range("G5").Value = 28.31746357381291
$ws.Range("H5").Value = 14.33036127678195
$ws.Range("I5").Value = 22.17614363388824
$ws.Range("J5").Value = 10.26027632857181
$ws.Range("L5").Value = 11.82091469798683
$ws.Range("O5").Value = 21.7254651464536
$ws.Range("B6").Value = 16.17058721929029
$ws.Range("C6").Value = 10.38827014165754
$ws.Range("D6").Value = 9.934879758172578
$ws.Range("F6").Value = 29.3804894056337
$ws.Range("G6").Value = 28.32521910566503
$ws.Range("H6").Value = 14.33328411104252
$ws.Range("I6").Value = 22.18297883448215
$ws.Range("J6").Value = 10.26172288850323
$ws.Range("L6").Value = 11.81957268011947
$ws.Range("O6").Value = 21.73074863150715
$ws.Range("B7").Value = 16.38093147656388
$ws.Range("C7").Value = 10.46184582183804
$ws.Range("D7").Value = 9.934072504538367
$ws.Range("F7").Value = 29.33792095083783
$ws.Range("G7").Value = 28.27234584470157
$ws.Range("H7").Value = 14.31322815243674
$ws.Range("I7").Value = 22.13597969757782
$ws.Range("J7").Value = 10.25177341108642
$ws.Range("L7").Value = 11.82899923308172
$ws.Range("O7").Value = 21.6945528980354
$ws.Range("B8").Value = 17.27718181351408
$ws.Range("C8").Value = 10.78022022019643
$ws.Range("D8").Value = 9.933340754596792
$ws.Range("F8").Value = 29.16726808690661
$ws.Range("G8").Value = 28.06286486058417
$ws.Range("H8").Value = 14.23045843609818
$ws.Range("I8").Value = 21.93957117778591
$ws.Range("J8").Value = 10.21012095623666
$ws.Range("L8").Value = 11.87332019618685
$ws.Range("O8").Value = 21.54666143796193
$ws.Range("B9").Value = 18.90646766213509
$ws.Range("C9").Value = 11.37736091608093
$ws.Range("D9").Value = 9.942003000028912
$ws.Range("F9").Value = 28.89478716910459
$ws.Range("G9").Value = 27.73883425755159
$ws.Range("H9").Value = 14.08874587194623
$ws.Range("I9").Value = 21.59377597712592
$ws.Range("J9").Value = 10.13648715087865
$ws.Range("L9").Value = 11.96967368481146
$ws.Range("O9").Value = 21.29935232872234
$ws.Range("B10").Value = 20.01560939014694
$ws.Range("C10").Value = 11.79581122396598
$ws.Range("D10").Value = 9.954448931545592
$ws.Range("F10").Value = 28.73282857780138
$ws.Range("G10").Value = 27.55461859474254
$ws.Range("H10").Value = 13.99718091851405
$ws.Range("I10").Value = 21.36358908751024
$ws.Range("J10").Value = 10.08724696529402
$ws.Range("L10").Value = 12.04589120814845
$ws.Range("O10").Value = 21.14385211581381
$ws.Range("B11").Value = 20.49948447181587
$ws.Range("C11").Value = 11.98107705173886
$ws.Range("D11").Value = 9.961413776196187
$ws.Range("F11").Value = 28.66751944849844
$ws.Range("G11").Value = 27.48270144416765
$ws.Range("H11").Value = 13.95825019152748
$ws.Range("I11").Value = 21.26402799881292
$ws.Range("J11").Value = 10.06589135947615
$ws.Range("L11").Value = 12.08166847542686
$ws.Range("O11").Value = 21.07883524846392
$ws.Range("B12").Value = 20.67963173985278
$ws.Range("C12").Value = 12.05044813962233
$ws.Range("D12").Value = 9.964236871818551
$ws.Range("F12").Value = 28.64399704905052
$ws.Range("G12").Value = 27.45719222300673
$ws.Range("H12").Value = 13.94389970589589
$ws.Range("I12").Value = 21.22706593278026
$ws.Range("J12").Value = 10.05795392788789
$ws.Range("L12").Value = 12.09536871816588
$ws.Range("O12").Value = 21.05504080160726
$ws.Range("B13").Value = 20.64097265296108
$ws.Range("C13").Value = 12.03554355271304
$ws.Range("D13").Value = 9.963620639021343
$ws.Range("F13").Value = 28.64900917174768
$ws.Range("G13").Value = 27.46260916607448
$ws.Range("H13").Value = 13.94697291335156
$ws.Range("I13").Value = 21.23499350299582
$ws.Range("J13").Value = 10.05965675829522
$ws.Range("L13").Value = 12.09241146835532
$ws.Range("O13").Value = 21.06012856985164
$ws.Range("B14").Value = 20.51436772259559
$ws.Range("C14").Value = 11.98680027642243
$ws.Range("D14").Value = 9.96164232333931
$ws.Range("F14").Value = 28.66555998987923
$ws.Range("G14").Value = 27.48056814299334
$ws.Range("H14").Value = 13.95706171609706
$ws.Range("I14").Value = 21.26097229276925
$ws.Range("J14").Value = 10.06523535026477
$ws.Range("L14").Value = 12.08279259260016
$ws.Range("O14").Value = 21.07686109081133
$ws.Range("B15").Value = 20.43641364746606
$ws.Range("C15").Value = 11.95683986722281
$ws.Range("D15").Value = 9.960454673201395
$ws.Range("F15").Value = 28.67585542239706
$ws.Range("G15").Value = 27.4917935277985
$ws.Range("H15").Value = 13.9632924217281
$ws.Range("I15").Value = 21.27698132639452
$ws.Range("J15").Value = 10.06867184541347
$ws.Range("L15").Value = 12.0769203651132
$ws.Range("O15").Value = 21.08721791659777
$ws.Range("B16").Value = 19.98356093794797
$ws.Range("C16").Value = 11.78359660942704
$ws.Range("D16").Value = 9.954019839525307
$ws.Range("F16").Value = 28.73726552814795
$ws.Range("G16").Value = 27.55955900890429
$ws.Range("H16").Value = 13.99977993803575
$ws.Range("I16").Value = 21.3701991985253
$ws.Range("J16").Value = 10.08866355635989
$ws.Range("L16").Value = 12.04357474817334
$ws.Range("O16").Value = 21.1482165283791
$ws.Range("B17").Value = 19.70036896652756
$ws.Range("C17").Value = 11.67597476260717
$ws.Range("D17").Value = 9.950404843868114
$ws.Range("F17").Value = 28.77708561716517
$ws.Range("G17").Value = 27.60418576323924
$ws.Range("H17").Value = 14.02286133076011
$ws.Range("I17").Value = 21.42870396263816
$ws.Range("J17").Value = 10.10119475065868
$ws.Range("L17").Value = 12.02339666283725
$ws.Range("O17").Value = 21.18710486275193
$ws.Range("B18").Value = 19.53554359950946
$ws.Range("C18").Value = 11.61359740107804
$ws.Range("D18").Value = 9.948448360688628
$ws.Range("F18").Value = 28.8007762636871
$ws.Range("G18").Value = 27.63097154086101
$ws.Range("H18").Value = 14.03639338918938
$ws.Range("I18").Value = 21.46283936721782
$ws.Range("J18").Value = 10.10850066723558
$ws.Range("L18").Value = 12.01189523945151
$ws.Range("O18").Value = 21.21001068554372
$ws.Range("B19").Value = 19.479406911664
$ws.Range("C19").Value = 11.59239741548068
$ws.Range("D19").Value = 9.947807065318285
$ws.Range("F19").Value = 28.80893253394606
$ws.Range("G19").Value = 27.64023222138434
$ws.Range("H19").Value = 14.04101911670138
$ws.Range("I19").Value = 21.47448037094129
$ws.Range("J19").Value = 10.11099122904743
$ws.Range("L19").Value = 12.00801920356784
$ws.Range("O19").Value = 21.21785855177156
$ws.Range("B20").Value = 19.73071696888752
$ws.Range("C20").Value = 11.68748096777369
$ws.Range("D20").Value = 9.950776972622641
$ws.Range("F20").Value = 28.77276518853226
$ws.Range("G20").Value = 27.5993193827416
$ws.Range("H20").Value = 14.02037775443907
$ws.Range("I20").Value = 21.42242584342061
$ws.Range("J20").Value = 10.0998506141784
$ws.Range("L20").Value = 12.02553389014641
$ws.Range("O20").Value = 21.1829093987502
$ws.Range("B21").Value = 20.55163921775041
$ws.Range("C21").Value = 12.00113906763186
$ws.Range("D21").Value = 9.962218378143138
$ws.Range("F21").Value = 28.66066576443028
$ws.Range("G21").Value = 27.47524624131164
$ws.Range("H21").Value = 13.95408775666743
$ws.Range("I21").Value = 21.25332163112606
$ws.Range("J21").Value = 10.06359273110354
$ws.Range("L21").Value = 12.08561381492254
$ws.Range("O21").Value = 21.07192389870715
$ws.Range("B22").Value = 21.07013254025042
$ws.Range("C22").Value = 12.20153512278702
$ws.Range("D22").Value = 9.97077729142552
$ws.Range("F22").Value = 28.59445028876026
$ws.Range("G22").Value = 27.40421412254862
$ws.Range("H22").Value = 13.91304686808867
$ws.Range("I22").Value = 21.14711240863789
$ws.Range("J22").Value = 10.04076705169825
$ws.Range("L22").Value = 12.12576318862581
$ws.Range("O22").Value = 21.00420476397528
$ws.Range("B23").Value = 20.79508525176957
$ws.Range("C23").Value = 12.09501715192839
$ws.Range("D23").Value = 9.966110897998389
$ws.Range("F23").Value = 28.62914403941784
$ws.Range("G23").Value = 27.44120018994259
$ws.Range("H23").Value = 13.93474214577899
$ws.Range("I23").Value = 21.20340430150084
$ws.Range("J23").Value = 10.05287006935797
$ws.Range("L23").Value = 12.10425615127473
$ws.Range("O23").Value = 21.03990593305932
$ws.Range("B24").Value = 19.71700290316744
$ws.Range("C24").Value = 11.6822805785957
$ws.Range("D24").Value = 9.950608353619
$ws.Range("F24").Value = 28.77471597044389
$ws.Range("G24").Value = 27.60151595697066
$ws.Range("H24").Value = 14.02149976238488
$ws.Range("I24").Value = 21.42526262279856
$ws.Range("J24").Value = 10.10045798185984
$ws.Range("L24").Value = 12.02456734009111
$ws.Range("O24").Value = 21.18480446003631
$ws.Range("B25").Value = 18.48052314673454
$ws.Range("C25").Value = 11.21912710440568
$ws.Range("D25").Value = 9.938585858058847
$ws.Range("F25").Value = 28.96180979742012
$ws.Range("G25").Value = 27.81710391558171
$ws.Range("H25").Value = 14.12487838350113
$ws.Range("I25").Value = 21.68312221121759
$ws.Range("J25").Value = 10.15555046965362
$ws.Range("L25").Value = 11.94262986385335
$ws.Range("O25").Value = 21.36166730342541
